$wb = $excel.ActiveWorkbook
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
$target = $wb.Worksheets.Item("Note Project")
$ws = $wb.Worksheets.Add($null, $target)
$ws.Name = "Note Tuan"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
